$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1 (style copied from existing header E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null

# Re-set values since PasteSpecial(Formats) should not overwrite values, but ensure correctness
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Boolean data values for rows 2-8
$data = @(
    @($false, $false, $false),
    @($false, $false, $false),
    @($true,  $true,  $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($true,  $true,  $true),
    @($false, $false, $false)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $data[$i][0]
    $ws.Cells.Item($row, 7).Value = $data[$i][1]
    $ws.Cells.Item($row, 8).Value = $data[$i][2]
}
